# Adds two new result columns (K: economic_losses, L: total_loss) to the
# Fairness results sheet, and rescales the existing "reward" column (J) which
# had accidentally been left multiplied by an extra factor of 1e9.
#
# Source data for this edit comes from a re-run of the underlying benchmark
# with the new "economic_losses"/"total_loss" objective tracked alongside the
# existing metrics, as described in the commit ("modified data files for new
# obj").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add labels for the two new columns -----------------
# Copy the formatting (bold font, border, centered alignment) from the
# existing "reward" header cell (J1) onto the two new header cells so they
# match the look of the rest of the header row.
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1:L1").PasteSpecial(-4122) | Out-Null
$ws.Range("K1").Value = "economic_losses"
$ws.Range("L1").Value = "total_loss"

# --- Data rows (rows 2-80): columns are [row, reward(J), economic_losses(K), total_loss(L)] ---
# reward(J) is $null for row 2, meaning that column is left unchanged there.
$data = @(
    @(2, $null, 98.01884882001667, 98.01884882001667),
    @(3, 55.70319825736978, 42.31565056264689, 42.31565056264689),
    @(4, 55.70319825736978, 42.31565056264689, 42.31565056264689),
    @(5, 55.70319825736978, 42.31565056264689, 42.31565056264689),
    @(6, 55.70319825736978, 42.31565056264689, 42.31565056264689),
    @(7, 55.70319825736978, 42.31565056264689, 42.31565056264689),
    @(8, 55.70319825736978, 42.31565056264689, 42.31565056264689),
    @(9, 55.70319825736978, 42.31565056264689, 42.31565056264689),
    @(10, 55.70319825736978, 42.31565056264689, 42.31565056264689),
    @(11, 55.70319825736978, 42.31565056264689, 42.31565056264689),
    @(12, 55.70319825736978, 42.31565056264689, 42.31565056264689),
    @(13, 55.70319825736978, 42.31565056264689, 42.31565056264689),
    @(14, 55.70319825736978, 42.31565056264689, 42.31565056264689),
    @(15, 55.70319825736978, 42.31565056264689, 42.31565056264689),
    @(16, 54.18103283177832, 42.26388329436435, 43.83781598823835),
    @(17, 54.65314602272955, 42.1394613833027, 43.36570279728712),
    @(18, 54.2791708377467, 42.23937634044267, 43.73967798226998),
    @(19, 54.44218411544126, 42.18334043924043, 43.57666470457541),
    @(20, 54.15603808066071, 42.27933732582964, 43.86281073935596),
    @(21, 54.25291823913898, 42.24518578260834, 43.7659305808777),
    @(22, 54.3344995345559, 42.22735353218471, 43.68434928546078),
    @(23, 54.59547871757277, 42.15278598465879, 43.42337010244391),
    @(24, 54.50438933482131, 42.14988680440002, 43.51445948519537),
    @(25, 54.97793171378946, 42.06043213573306, 43.04091710622722),
    @(26, 54.97793171378946, 42.06043213573306, 43.04091710622722),
    @(27, 54.22763738523792, 42.25086110861932, 43.79121143477876),
    @(28, 54.39379052282607, 42.21474421636126, 43.62505829719061),
    @(29, 53.25983419956586, 43.15547325137542, 44.75901462045082),
    @(30, 52.67248548369619, 43.60677653232785, 45.34636333632049),
    @(31, 52.05110021217882, 44.2169154260073, 45.96774860783786),
    @(32, 52.89963913036214, 43.38711215166814, 45.11920968965454),
    @(33, 51.22666447670407, 44.76219647125619, 46.79218434331261),
    @(34, 51.58337454851122, 44.65907588686757, 46.43547427150546),
    @(35, 51.67419974019302, 44.56199839316476, 46.34464907982365),
    @(36, 50.87041424058946, 44.72079213772634, 47.14843457942722),
    @(37, 52.46036875639086, 43.82530348219947, 45.55848006362582),
    @(38, 52.04464303984525, 42.40776421902063, 45.97420578017142),
    @(39, 51.86002668040393, 44.36753427181317, 46.15882213961275),
    @(40, 53.25983419956586, 43.15547325137542, 44.75901462045082),
    @(41, 52.25099859983963, 44.02144141367193, 45.76785022017705),
    @(42, 49.90841011483423, 45.19993003497599, 48.11043870518245),
    @(43, 50.53644164308379, 44.1923695243806, 47.48240717693287),
    @(44, 49.6608546797615, 45.63168923904469, 48.35799414025518),
    @(45, 51.44089354248672, 44.1744375025549, 46.57795527752996),
    @(46, 50.03660472535678, 45.01953938531658, 47.9822440946599),
    @(47, 50.92804874606259, 44.19494399001952, 47.09080007395409),
    @(48, 49.70459018202234, 45.62217869007696, 48.31425863799433),
    @(49, 49.63351991782667, 45.73983665058844, 48.38532890219001),
    @(50, 50.1587078598961, 44.86221450242451, 47.86014096012057),
    @(51, 49.91600681549107, 45.17174735259648, 48.1028420045256),
    @(52, 51.44089354248672, 44.1744375025549, 46.57795527752996),
    @(53, 50.41881298683839, 44.6005429621354, 47.60003583317829),
    @(54, 49.72698668777207, 45.62246959087071, 48.29186213224462),
    @(55, 47.49864545123842, 46.05317106863232, 50.52020336877825),
    @(56, 47.55073130381823, 46.04684070148503, 50.46811751619845),
    @(57, 48.30083684226291, 45.45871670435888, 49.71801197775377),
    @(58, 47.4260251145131, 46.06255206619931, 50.59282370550358),
    @(59, 47.63775868147027, 45.94382836988303, 50.38109013854641),
    @(60, 47.70227681507964, 45.78735819914484, 50.31657200493703),
    @(61, 47.40602342066192, 46.05012817623154, 50.61282539935475),
    @(62, 47.57937403320069, 46.02017060516717, 50.43947478681599),
    @(63, 47.31349803178252, 46.08333682634265, 50.70535078823416),
    @(64, 47.45280703138995, 46.05419475438281, 50.56604178862672),
    @(65, 48.30083684226291, 45.45871670435888, 49.71801197775377),
    @(66, 47.81225373698629, 45.69988424781952, 50.20659508303039),
    @(67, 47.3567573019763, 46.07191552235265, 50.66209151804037),
    @(68, 44.33146335280103, 47.88070679301862, 53.68738546721563),
    @(69, 44.02215382559327, 48.01161238253772, 53.99669499442341),
    @(70, 44.53366724335636, 45.97179671410257, 53.48518157666032),
    @(71, 43.79618340953561, 46.82451145785758, 54.22266541048107),
    @(72, 43.78129342551886, 47.11764746974928, 54.23755539449782),
    @(73, 44.58323562182331, 47.63228553761971, 53.43561319819337),
    @(74, 44.96272648419868, 47.43396707994586, 53.056122335818),
    @(75, 44.73443309045485, 47.55449282374046, 53.28441572956183),
    @(76, 44.96272648419868, 47.43396707994586, 53.056122335818),
    @(77, 43.84866843266285, 46.92434745507822, 54.17018038735382),
    @(78, 43.67904302031744, 47.01075966206562, 54.33980579969924),
    @(79, 44.45552837612134, 47.75652134205526, 53.56332044389534),
    @(80, 44.20984720540562, 48.00485863202475, 53.80900161461106)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $j = $entry[1]
    $k = $entry[2]
    $l = $entry[3]

    if ($null -ne $j) {
        $ws.Cells.Item($r, 10).Value = $j   # column J = 10
    }
    $ws.Cells.Item($r, 11).Value = $k       # column K = 11
    $ws.Cells.Item($r, 12).Value = $l       # column L = 12
}
